# Auto-generated edit script applying the profit-sheet value updates
# described by the commit diff (scheduled runner refresh of price data).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3580.7
$ws.Range("I74").Value = 3160.6
$ws.Range("J74").Value = 4000.8
$ws.Range("K74").Value = 3160.6
$ws.Range("L74").Value = 4000.8
$ws.Range("M74").Value = -2224.6
$ws.Range("N74").Value = -5872.8
$ws.Range("H77").Value = 3580.7
$ws.Range("I77").Value = 3160.6
$ws.Range("J77").Value = 4000.8
$ws.Range("K77").Value = 15803
$ws.Range("L77").Value = 20004
$ws.Range("M77").Value = -11123
$ws.Range("N77").Value = -29364
$ws.Range("H112").Value = 13464.632
$ws.Range("I112").Value = 784.6667
$ws.Range("J112").Value = 15842.125
$ws.Range("K112").Value = 2354.0001
$ws.Range("L112").Value = 47526.375
$ws.Range("M112").Value = -1246.0001
$ws.Range("N112").Value = -49742.375
$ws.Range("H129").Value = 735.6667
$ws.Range("J129").Value = 914.7273
$ws.Range("L129").Value = 2744.1819
$ws.Range("N129").Value = -12744.1819
$ws.Range("H138").Value = 2034.1957
$ws.Range("I138").Value = 1327.0526
$ws.Range("J138").Value = 2218.2466
$ws.Range("K138").Value = 3981.1578
$ws.Range("L138").Value = 6654.739799999999
$ws.Range("M138").Value = 1158.8422
$ws.Range("N138").Value = -16934.7398

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 273.14285
$ws.Range("I5").Value = 175.5
$ws.Range("K5").Value = 175.5
$ws.Range("M5").Value = -63.5
$ws.Range("H23").Value = 86673
$ws.Range("I23").Value = 80006
$ws.Range("K23").Value = 80006
$ws.Range("M23").Value = -79747
$ws.Range("H37").Value = 28000
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").Value = $null
$ws.Range("H61").Value = 1367.3334
$ws.Range("I61").Value = 1135.5834
$ws.Range("K61").Value = 1135.5834
$ws.Range("M61").Value = -923.5834
$ws.Range("H63").Value = 2696.2693
$ws.Range("I63").Value = 2333.925
$ws.Range("J63").Value = 3904.0833
$ws.Range("K63").Value = 2333.925
$ws.Range("L63").Value = 3904.0833
$ws.Range("M63").Value = -1647.925
$ws.Range("N63").Value = -5276.0833
$ws.Range("H66").Value = 2696.2693
$ws.Range("I66").Value = 2333.925
$ws.Range("J66").Value = 3904.0833
$ws.Range("K66").Value = 11669.625
$ws.Range("L66").Value = 19520.4165
$ws.Range("M66").Value = -8237.625
$ws.Range("N66").Value = -26384.4165
$ws.Range("H80").Value = 38000
$ws.Range("J80").Value = 38000
$ws.Range("L80").Value = 38000
$ws.Range("N80").Value = -39996
$ws.Range("H83").Value = 38000
$ws.Range("J83").Value = 38000
$ws.Range("L83").Value = 114000
$ws.Range("N83").Value = -123984
$ws.Range("H122").Value = 3142.8667
$ws.Range("I122").Value = 3101.3845
$ws.Range("J122").Value = 3412.5
$ws.Range("K122").Value = 9304.1535
$ws.Range("L122").Value = 10237.5
$ws.Range("M122").Value = -6854.1535
$ws.Range("N122").Value = -15137.5
$ws.Range("H132").Value = 3584.0322
$ws.Range("I132").Value = 2980.2856
$ws.Range("K132").Value = 8940.856800000001
$ws.Range("M132").Value = -6410.856800000001
$ws.Range("H136").Value = 1367.3334
$ws.Range("I136").Value = 1135.5834
$ws.Range("K136").Value = 3406.7502
$ws.Range("M136").Value = -856.7501999999999
$ws.Range("H138").Value = 56410
$ws.Range("J138").Value = 56410
$ws.Range("L138").Value = 56410
$ws.Range("N138").Value = -66690

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 273.14285
$ws.Range("I4").Value = 175.5
$ws.Range("K4").Value = 175.5
$ws.Range("M4").Value = -60.5
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = -127
$ws.Range("H94").Value = 15625700
$ws.Range("I94").Value = 22727928
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 22727928
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = -22727477
$ws.Range("N94").Value = -1702
$ws.Range("H134").Value = 1194.8334
$ws.Range("I134").Value = 1025.4375
$ws.Range("J134").Value = 2550
$ws.Range("K134").Value = 3076.3125
$ws.Range("L134").Value = 7650
$ws.Range("M134").Value = -541.3125
$ws.Range("N134").Value = -12720

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1036.3572
$ws.Range("I122").Value = 983.1111
$ws.Range("K122").Value = 2949.3333
$ws.Range("M122").Value = -499.3332999999998
$ws.Range("H134").Value = 2127.2104
$ws.Range("I134").Value = 1977.3125
$ws.Range("K134").Value = 5931.9375
$ws.Range("M134").Value = -3396.9375
$ws.Range("H141").Value = 618862.9
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 618862.9
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 618862.9
$ws.Range("M141").Value = $null
$ws.Range("N141").Value = -629222.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 349.7143
$ws.Range("I13").Value = 89.8
$ws.Range("J13").Value = 999.5
$ws.Range("K13").Value = 269.4
$ws.Range("L13").Value = 2998.5
$ws.Range("M13").Value = -101.4
$ws.Range("N13").Value = -3334.5
$ws.Range("H50").Value = 284.0909
$ws.Range("I50").Value = 148.75
$ws.Range("J50").Value = 361.42856
$ws.Range("K50").Value = 446.25
$ws.Range("L50").Value = 1084.28568
$ws.Range("M50").Value = 34.75
$ws.Range("N50").Value = -2046.28568
$ws.Range("H53").Value = 284.0909
$ws.Range("I53").Value = 148.75
$ws.Range("J53").Value = 361.42856
$ws.Range("K53").Value = 446.25
$ws.Range("L53").Value = 1084.28568
$ws.Range("M53").Value = 34.75
$ws.Range("N53").Value = -2046.28568
$ws.Range("H63").Value = 13381.583
$ws.Range("H66").Value = 13381.583
$ws.Range("H74").Value = 5250
$ws.Range("J74").Value = 5250
$ws.Range("L74").Value = 15750
$ws.Range("N74").Value = -17872
$ws.Range("H77").Value = 5250
$ws.Range("J77").Value = 5250
$ws.Range("L77").Value = 47250
$ws.Range("N77").Value = -57858
$ws.Range("H81").Value = 3411.1875
$ws.Range("J81").Value = 3621.5386
$ws.Range("L81").Value = 10864.6158
$ws.Range("N81").Value = -13110.6158
$ws.Range("H84").Value = 3411.1875
$ws.Range("J84").Value = 3621.5386
$ws.Range("L84").Value = 32593.8474
$ws.Range("N84").Value = -43825.8474
$ws.Range("H113").Value = 735.2
$ws.Range("J113").Value = 741.2083
$ws.Range("L113").Value = 2223.6249
$ws.Range("N113").Value = -6563.6249
$ws.Range("H122").Value = 1241.0416
$ws.Range("I122").Value = 652.7778
$ws.Range("K122").Value = 5875.000199999999
$ws.Range("M122").Value = -3425.000199999999
$ws.Range("H132").Value = 1129
$ws.Range("I132").Value = 843
$ws.Range("K132").Value = 7587
$ws.Range("M132").Value = -5057

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 271.9
$ws.Range("J107").Value = 324.5
$ws.Range("L107").Value = 324.5
$ws.Range("N107").Value = -4164.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 7000
$ws.Range("J15").Value = 7000
$ws.Range("L15").Value = 7000
$ws.Range("N15").Value = -7576
$ws.Range("H81").Value = 1207.9
$ws.Range("I81").Value = 1207.9
$ws.Range("K81").Value = 2415.8
$ws.Range("M81").Value = -1354.8
$ws.Range("H84").Value = 1207.9
$ws.Range("I84").Value = 1207.9
$ws.Range("K84").Value = 12079
$ws.Range("M84").Value = -6775
$ws.Range("H141").Value = 50300
$ws.Range("J141").Value = 50300
$ws.Range("L141").Value = 50300
$ws.Range("N141").Value = -60660
